$d = $word.ActiveDocument
$d.Save()
Write-Host "saved"
